# This edit re-shuffles the per-row data (Fecha, Calidad, Volumen, Precio
# minimo/maximo/promedio, Unidad de comercializacion, Origen, Precio $/Kg,
# Kg/unidad) across rows 2-18 of the sheet, while leaving the descriptive
# columns (Mercado ID, Mercado, Region, Codreg, Tipo, Producto ID, Producto,
# Categoria ID, Categoria, Variedad) untouched. Row 4 keeps its own data.
#
# Mapping: new row N gets the old data that used to live in row Map[N].
$map = @{
    2  = 15
    3  = 13
    4  = 4
    5  = 9
    6  = 12
    7  = 8
    8  = 7
    9  = 11
    10 = 6
    11 = 5
    12 = 17
    13 = 18
    14 = 3
    15 = 16
    16 = 2
    17 = 10
    18 = 14
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together with the row's "record".
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# 1. Snapshot the current ("before") values for every affected column/row.
$snapshot = @{}
foreach ($col in $cols) {
    $snapshot[$col] = @{}
    for ($r = 2; $r -le 18; $r++) {
        $snapshot[$col][$r] = $ws.Range("$col$r").Value()
    }
}

# 2. Write back the values according to the permutation map.
foreach ($col in $cols) {
    for ($r = 2; $r -le 18; $r++) {
        $src = $map[$r]
        $ws.Range("$col$r").Value = $snapshot[$col][$src]
    }
}
